$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Change selection to I6
$ws.Range("I6").Select()

# 2. Change column J width
$ws.Columns.Item(10).ColumnWidth = 0.42578125

# 3. Update text of C6
$ws.Range("C6").Value = "The columns are set as [1, 2, 4, 8, 16, 32, 16, 21, 30, .5]"

# 4. Page setup scale
$ws.PageSetup.Zoom = $false
$ws.PageSetup.Orientation = 2
$ws.PageSetup.PrintScale = 86

Write-Host "done"
